$wb = $excel.ActiveWorkbook
$wsButtons   = $wb.Worksheets.Item("Buttons")
$wsLabels    = $wb.Worksheets.Item("Labels")
$wsTextBoxes = $wb.Worksheets.Item("TextBoxes")

# ---------------------------------------------------------------------------
# TextBoxes sheet: extend the header/template rows with the same grid-
# placement columns (rowStart/rowEnd/columnStart/c9olumnEnd/placement) that
# the Buttons and Labels sheets already use, pushing the old "textBoxType"
# column out to column I. The old "height"/"width" columns are dropped.
# ---------------------------------------------------------------------------

# 1) Preserve the old F1 header ("textBoxType") by moving it out to I1 first
#    (copies both value and style) before it gets overwritten below.
$wsTextBoxes.Range("F1").Copy($wsTextBoxes.Range("I1"))

# 2) Overwrite D1:H1 with the rowStart/rowEnd/columnStart/c9olumnEnd/placement
#    headers, copying both value and style straight from the Buttons sheet
#    (which already has this exact header layout).
$wsButtons.Range("D1:H1").Copy($wsTextBoxes.Range("D1"))

# 3) Extend the blank style-template rows 2 and 3 across the new columns.
#    Row 2 uses the "s=2" blank style (bordered), row 3 the "s=1" blank
#    style - reuse the already-blank B2/B3 cells on this sheet as the
#    style source so no values get dragged along.
$wsTextBoxes.Range("B2").Copy($wsTextBoxes.Range("D2"))
$wsTextBoxes.Range("B2").Copy($wsTextBoxes.Range("E2"))
$wsTextBoxes.Range("B2").Copy($wsTextBoxes.Range("G2"))
$wsTextBoxes.Range("B2").Copy($wsTextBoxes.Range("H2"))
$wsTextBoxes.Range("B2").Copy($wsTextBoxes.Range("I2"))

$wsTextBoxes.Range("B3").Copy($wsTextBoxes.Range("D3"))
$wsTextBoxes.Range("B3").Copy($wsTextBoxes.Range("E3"))
$wsTextBoxes.Range("B3").Copy($wsTextBoxes.Range("G3"))
$wsTextBoxes.Range("B3").Copy($wsTextBoxes.Range("H3"))
$wsTextBoxes.Range("B3").Copy($wsTextBoxes.Range("I3"))

# 4) Column widths: drop the old custom widths for columns E/F (old
#    "width"/"textBoxType" columns) and give column I (new "textBoxType")
#    the width the old column F used to have.
$wsTextBoxes.Columns.Item(5).ColumnWidth = $wsTextBoxes.Columns.Item(1).ColumnWidth
$wsTextBoxes.Columns.Item(5).NumberFormat = "General"
$wsTextBoxes.Columns.Item(9).ColumnWidth = 11.08984375
$wsTextBoxes.Columns.Item(9).SetAttr("bestFit", "1")

# ---------------------------------------------------------------------------
# Active tab / selection: the TextBoxes sheet becomes the active tab with a
# new zoom level and selection; Buttons loses its old tabSelected flag.
# ---------------------------------------------------------------------------
$wsTextBoxes.Activate()
$excel.ActiveWindow.Zoom = 92
$wsTextBoxes.Range("G9").Select()
